$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '62.864.70'
$ws.Range("E2").Value = '  -1.47%  '
$ws.Range("D3").Value = '3.171.39'
$ws.Range("E3").Value = '  -4.47%  '
$ws.Range("E4").Value = '  +0.06%  '
$ws.Range("D5").Value = '590.05'
$ws.Range("E5").Value = '  -2.35%  '
$ws.Range("D6").Value = '134.26'
$ws.Range("E6").Value = '  -5.90%  '
$ws.Range("D8").Value = '3.170.64'
$ws.Range("E8").Value = '  -4.39%  '
$ws.Range("E9").Value = '  -0.66%  '
$ws.Range("E10").Value = '  -6.57%  '
$ws.Range("E11").Value = '  -5.71%  '
$ws.Range("E12").Value = '  -3.84%  '
$ws.Range("E13").Value = '  -5.24%  '
$ws.Range("D14").Value = '34.71'
$ws.Range("E14").Value = '  -0.98%  '
$ws.Range("D15").Value = '3.692.36'
$ws.Range("E15").Value = '  -4.53%  '
$ws.Range("E16").Value = '  -1.08%  '
$ws.Range("D17").Value = '3.165.34'
$ws.Range("E17").Value = '  -4.74%  '
$ws.Range("D18").Value = '62.851.43'
$ws.Range("E18").Value = '  -1.65%  '
$ws.Range("D19").Value = '6.54'
$ws.Range("E19").Value = '  -4.72%  '
$ws.Range("D20").Value = '459.31'
$ws.Range("E20").Value = '  -4.65%  '
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '14.00'
$ws.Range("E21").Value = '  -0.68%  '
$ws.Range("D22").Value = '0.692'
$ws.Range("E22").Value = '  -6.33%  '
$ws.Range("E23").Value = '  -4.60%  '
$ws.Range("D24").Value = '13.35'
$ws.Range("E24").Value = '  -4.32%  '
$ws.Range("D25").Value = '82.78'
$ws.Range("E25").Value = '  -2.61%  '
$ws.Range("E26").Value = '  +0.00%  '
$ws.Range("E27").Value = '  +0.04%  '
$ws.Range("E28").Value = '  -4.17%  '
$ws.Range("E29").Value = '  -6.73%  '
$ws.Range("E30").Value = '  -5.67%  '
$ws.Range("E31").Value = '  -6.07%  '
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '27.10'
$ws.Range("E32").Value = '  -6.21%  '
$ws.Range("E33").Value = '  -3.72%  '
$ws.Range("D34").Value = '2.36'
$ws.Range("E34").Value = '  -6.73%  '
$ws.Range("E35").Value = '  -6.24%  '
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '5.80'
$ws.Range("E36").Value = '  -4.82%  '
$ws.Range("D37").Value = '51.29'
$ws.Range("E37").Value = '  -2.17%  '
$ws.Range("D38").Value = '0.0₃0702'
$ws.Range("E38").Value = '  -5.76%  '
$ws.Range("E39").Value = '  -3.28%  '
$ws.Range("D40").Value = '404.61'
$ws.Range("E40").Value = '  -6.59%  '
$ws.Range("D41").Value = '8.06'
$ws.Range("E41").Value = '  -3.66%  '
$ws.Range("E42").Value = '  -5.33%  '
$ws.Range("E43").Value = '  -5.14%  '
$ws.Range("D44").Value = '2.791.18'
$ws.Range("E44").Value = '  -10.92%  '
$ws.Range("E45").Value = '  -6.64%  '
$ws.Range("E46").Value = '  -0.01%  '
$ws.Range("E47").Value = '  -6.34%  '
$ws.Range("D48").Value = '124.76'
$ws.Range("E48").Value = '  +0.19%  '
$ws.Range("D49").Value = '25.17'
$ws.Range("E49").Value = '  -4.67%  '
$ws.Range("D50").Value = '34.35'
$ws.Range("E50").Value = '  -6.29%  '
